$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Mean" column (column 52 / AZ),
# which shifts "Mean" and its data one column to the right (-> BA) and
# leaves a blank column 52 (AZ) ready for the new "Run 50" data.
$ws.Columns.Item(52).Insert()

# New header for the inserted column.
$ws.Range("AZ1").Value = "Run 50"

# Fill the new "Run 50" column and refresh the (now shifted) "Mean" column
# for every data row (2-14).
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 52).Value = 823322052.5601799
    $ws.Cells.Item($r, 53).Value = 1471355866.644745
}
